# Add new row 61 with the 2025-02-10 18:43:34 resale-number update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (A-D) ---
# Force these cells to be stored as plain text (otherwise Excel's smart entry
# parsing would turn "2025-02-10" / "18:43:34" into date/time serials and
# "06" into the number 6), then restore the default (unstyled) cell
# formatting by copying it from the row directly above, which already uses
# the workbook's default style.
$ws.Range("A61:D61").NumberFormat = "@"
$ws.Range("A61").Value = "2025-02-10"
$ws.Range("B61").Value = "18:43:34"
$ws.Range("C61").Value = "Monday"
$ws.Range("D61").Value = "06"
$ws.Range("A61:D61").Style = $ws.Range("A60:D60").Style

# --- Numeric columns (E-T) ---
$ws.Range("E61").Value = 127536
$ws.Range("F61").Value = 141955
$ws.Range("G61").Value = 169299
$ws.Range("H61").Value = 158477
$ws.Range("I61").Value = -1
$ws.Range("J61").Value = 144277
$ws.Range("K61").Value = -1
$ws.Range("L61").Value = -1
$ws.Range("M61").Value = 191707
$ws.Range("N61").Value = 115027
$ws.Range("O61").Value = 44922
$ws.Range("P61").Value = 28484
$ws.Range("Q61").Value = 64757
$ws.Range("R61").Value = -1
$ws.Range("S61").Value = 43235
$ws.Range("T61").Value = -1
